$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 26.86490966666667
$ws.Range("H2").Value = 80.594729
$ws.Range("I2").Value = 0.1447302967754861
$ws.Range("J2").Value = 0.1447302967754861
$ws.Range("M2").Value = 1.442875
$ws.Range("N2").Value = 4.328625
$ws.Range("O2").Value = 0.02047893724893121
$ws.Range("P2").Value = 0.02047893724893121
$ws.Range("Q2").Value = 38.76270653529166
$ws.Range("R2").Value = 348.864358817625
$ws.Range("S2").Value = 0.002963922665684371
$ws.Range("T2").Value = 0.002963922665684371
$ws.Range("G3").Value = 26.86490966666667
$ws.Range("H3").Value = 80.594729
$ws.Range("I3").Value = 0.1447302967754861
$ws.Range("J3").Value = 0.1447302967754861
$ws.Range("O3").Value = 0.1473796107804731
$ws.Range("P3").Value = 0.1473796107804731
$ws.Range("Q3").Value = 278.9613802965844
$ws.Range("R3").Value = 2510.652422669259
$ws.Range("S3").Value = 0.0213302948069135
$ws.Range("T3").Value = 0.0213302948069135
$ws.Range("G4").Value = 26.86490966666667
$ws.Range("H4").Value = 80.594729
$ws.Range("I4").Value = 0.1447302967754861
$ws.Range("J4").Value = 0.1447302967754861
$ws.Range("M4").Value = 27.934719
$ws.Range("N4").Value = 83.804157
$ws.Range("O4").Value = 0.3964815784233052
$ws.Range("P4").Value = 0.3964815784233051
$ws.Range("Q4").Value = 750.463702498717
$ws.Range("R4").Value = 6754.173322488454
$ws.Range("S4").Value = 0.05738289651121813
$ws.Range("T4").Value = 0.05738289651121813
$ws.Range("G5").Value = 26.86490966666667
$ws.Range("H5").Value = 80.594729
$ws.Range("I5").Value = 0.1447302967754861
$ws.Range("J5").Value = 0.1447302967754861
$ws.Range("M5").Value = 30.695086
$ws.Range("N5").Value = 92.085258
$ws.Range("O5").Value = 0.4356598735472906
$ws.Range("P5").Value = 0.4356598735472905
$ws.Range("Q5").Value = 824.6207126005646
$ws.Range("R5").Value = 7421.586413405082
$ws.Range("S5").Value = 0.06305318279167012
$ws.Range("T5").Value = 0.06305318279167012
$ws.Range("I6").Value = 0.4077186109324291
$ws.Range("J6").Value = 0.4077186109324292
$ws.Range("M6").Value = 1.442875
$ws.Range("N6").Value = 4.328625
$ws.Range("O6").Value = 0.02047893724893121
$ws.Range("P6").Value = 0.02047893724893121
$ws.Range("Q6").Value = 109.1981237975833
$ws.Range("R6").Value = 982.78311417825
$ws.Range("S6").Value = 0.008349643848506614
$ws.Range("T6").Value = 0.008349643848506614
$ws.Range("I7").Value = 0.4077186109324291
$ws.Range("J7").Value = 0.4077186109324292
$ws.Range("O7").Value = 0.1473796107804731
$ws.Range("P7").Value = 0.1473796107804731
$ws.Range("Q7").Value = 785.8599685921528
$ws.Range("R7").Value = 7072.739717329375
$ws.Range("S7").Value = 0.06008941018717655
$ws.Range("T7").Value = 0.06008941018717654
$ws.Range("I8").Value = 0.4077186109324291
$ws.Range("J8").Value = 0.4077186109324292
$ws.Range("M8").Value = 27.934719
$ws.Range("N8").Value = 83.804157
$ws.Range("O8").Value = 0.3964815784233052
$ws.Range("P8").Value = 0.3964815784233051
$ws.Range("Q8").Value = 2114.125550454962
$ws.Range("R8").Value = 19027.12995409466
$ws.Range("S8").Value = 0.1616529184150469
$ws.Range("T8").Value = 0.1616529184150469
$ws.Range("I9").Value = 0.4077186109324291
$ws.Range("J9").Value = 0.4077186109324292
$ws.Range("M9").Value = 30.695086
$ws.Range("N9").Value = 92.085258
$ws.Range("O9").Value = 0.4356598735472906
$ws.Range("P9").Value = 0.4356598735472905
$ws.Range("Q9").Value = 2323.032695836761
$ws.Range("R9").Value = 20907.29426253085
$ws.Range("S9").Value = 0.177626638481699
$ws.Range("T9").Value = 0.177626638481699
$ws.Range("G10").Value = 14.45399366666666
$ws.Range("H10").Value = 43.36198099999999
$ws.Range("I10").Value = 0.07786852138807973
$ws.Range("J10").Value = 0.07786852138807973
$ws.Range("M10").Value = 1.442875
$ws.Range("N10").Value = 4.328625
$ws.Range("O10").Value = 0.02047893724893121
$ws.Range("P10").Value = 0.02047893724893121
$ws.Range("Q10").Value = 20.85530611179166
$ws.Range("R10").Value = 187.697755006125
$ws.Range("S10").Value = 0.001594664563173542
$ws.Range("T10").Value = 0.001594664563173542
$ws.Range("G11").Value = 14.45399366666666
$ws.Range("H11").Value = 43.36198099999999
$ws.Range("I11").Value = 0.07786852138807973
$ws.Range("J11").Value = 0.07786852138807973
$ws.Range("O11").Value = 0.1473796107804731
$ws.Range("P11").Value = 0.1473796107804731
$ws.Range("Q11").Value = 150.0882033135723
$ws.Range("R11").Value = 1350.793829822151
$ws.Range("S11").Value = 0.01147623237422614
$ws.Range("T11").Value = 0.01147623237422613
$ws.Range("G12").Value = 14.45399366666666
$ws.Range("H12").Value = 43.36198099999999
$ws.Range("I12").Value = 0.07786852138807973
$ws.Range("J12").Value = 0.07786852138807973
$ws.Range("M12").Value = 27.934719
$ws.Range("N12").Value = 83.804157
$ws.Range("O12").Value = 0.3964815784233052
$ws.Range("P12").Value = 0.3964815784233051
$ws.Range("Q12").Value = 403.768251506113
$ws.Range("R12").Value = 3633.914263555017
$ws.Range("S12").Value = 0.03087343426943475
$ws.Range("T12").Value = 0.03087343426943475
$ws.Range("G13").Value = 14.45399366666666
$ws.Range("H13").Value = 43.36198099999999
$ws.Range("I13").Value = 0.07786852138807973
$ws.Range("J13").Value = 0.07786852138807973
$ws.Range("M13").Value = 30.695086
$ws.Range("N13").Value = 92.085258
$ws.Range("O13").Value = 0.4356598735472906
$ws.Range("P13").Value = 0.4356598735472905
$ws.Range("Q13").Value = 443.6665786417886
$ws.Range("R13").Value = 3992.999207776097
$ws.Range("S13").Value = 0.03392419018124531
$ws.Range("T13").Value = 0.0339241901812453
$ws.Range("G14").Value = 68.62066266666666
$ws.Range("H14").Value = 205.861988
$ws.Range("I14").Value = 0.369682570904005
$ws.Range("J14").Value = 0.369682570904005
$ws.Range("M14").Value = 1.442875
$ws.Range("N14").Value = 4.328625
$ws.Range("O14").Value = 0.02047893724893121
$ws.Range("P14").Value = 0.02047893724893121
$ws.Range("Q14").Value = 99.01103864516665
$ws.Range("R14").Value = 891.0993478065
$ws.Range("S14").Value = 0.00757070617156668
$ws.Range("T14").Value = 0.00757070617156668
$ws.Range("G15").Value = 68.62066266666666
$ws.Range("H15").Value = 205.861988
$ws.Range("I15").Value = 0.369682570904005
$ws.Range("J15").Value = 0.369682570904005
$ws.Range("O15").Value = 0.1473796107804731
$ws.Range("P15").Value = 0.1473796107804731
$ws.Range("Q15").Value = 712.5471483759053
$ws.Range("R15").Value = 6412.924335383148
$ws.Range("S15").Value = 0.0544836734121569
$ws.Range("T15").Value = 0.05448367341215689
$ws.Range("G16").Value = 68.62066266666666
$ws.Range("H16").Value = 205.861988
$ws.Range("I16").Value = 0.369682570904005
$ws.Range("J16").Value = 0.369682570904005
$ws.Range("M16").Value = 27.934719
$ws.Range("N16").Value = 83.804157
$ws.Range("O16").Value = 0.3964815784233052
$ws.Range("P16").Value = 0.3964815784233051
$ws.Range("Q16").Value = 1916.898929187124
$ws.Range("R16").Value = 17252.09036268412
$ws.Range("S16").Value = 0.1465723292276053
$ws.Range("T16").Value = 0.1465723292276053
$ws.Range("G17").Value = 68.62066266666666
$ws.Range("H17").Value = 205.861988
$ws.Range("I17").Value = 0.369682570904005
$ws.Range("J17").Value = 0.369682570904005
$ws.Range("M17").Value = 30.695086
$ws.Range("N17").Value = 92.085258
$ws.Range("O17").Value = 0.4356598735472906
$ws.Range("P17").Value = 0.4356598735472905
$ws.Range("Q17").Value = 2106.317141930323
$ws.Range("R17").Value = 18956.8542773729
$ws.Range("S17").Value = 0.1610558620926761
$ws.Range("T17").Value = 0.1610558620926761
